$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "R") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the target paragraph containing only 'R'"
}

$startPos = $target.Range.Start

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:rFonts w:ascii="-apple-system" w:hAnsi="-apple-system" w:eastAsia="-apple-system" w:cs="-apple-system" w:asciiTheme="minorAscii" w:hAnsiTheme="minorAscii" w:eastAsiaTheme="minorAscii" w:cstheme="minorAscii"/><w:color w:val="24292E"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="-apple-system" w:hAnsi="-apple-system" w:eastAsia="-apple-system" w:cs="-apple-system"/><w:noProof w:val="0"/><w:color w:val="24292E"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">to design the output </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="-apple-system" w:hAnsi="-apple-system" w:eastAsia="-apple-system" w:cs="-apple-system"/><w:noProof w:val="0"/><w:color w:val="24292E"/><w:lang w:val="en-GB"/></w:rPr><w:t>impandance</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="-apple-system" w:hAnsi="-apple-system" w:eastAsia="-apple-system" w:cs="-apple-system"/><w:noProof w:val="0"/><w:color w:val="24292E"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> to be </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="-apple-system" w:hAnsi="-apple-system" w:eastAsia="-apple-system" w:cs="-apple-system"/><w:b w:val="1"/><w:bCs w:val="1"/><w:i w:val="1"/><w:iCs w:val="1"/><w:noProof w:val="0"/><w:color w:val="24292E"/><w:lang w:val="en-GB"/></w:rPr><w:t>50ohms</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="-apple-system" w:hAnsi="-apple-system" w:eastAsia="-apple-system" w:cs="-apple-system" w:asciiTheme="minorAscii" w:hAnsiTheme="minorAscii" w:eastAsiaTheme="minorAscii" w:cstheme="minorAscii"/><w:color w:val="24292E"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="-apple-system" w:hAnsi="-apple-system" w:eastAsia="-apple-system" w:cs="-apple-system"/><w:noProof w:val="0"/><w:color w:val="24292E"/><w:lang w:val="en-GB"/></w:rPr><w:t>May need a bandpass near output</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:eastAsia="Arial" w:cs="Arial" w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi" w:eastAsiaTheme="minorBidi" w:cstheme="minorBidi"/><w:noProof w:val="0"/><w:color w:val="24292E"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="-apple-system" w:hAnsi="-apple-system" w:eastAsia="-apple-system" w:cs="-apple-system"/><w:noProof w:val="0"/><w:color w:val="24292E"/><w:lang w:val="en-GB"/></w:rPr><w:t>Think how</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.Range.InsertXML($xml)

# InsertXML normalizes paragraph-mark / run sizes that equal the document's
# default run size (22 half-points / 11pt) by omitting them. The target
# formatting needs them explicit, so re-assert them now that the paragraphs
# exist. This forces Word to persist sz/szCs=22 on both the paragraph mark
# and the runs of the first two new paragraphs.

$p1 = $d.Paragraphs(1)
$found1 = $null
$found2 = $null
$found3 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Start -ge $startPos) {
        $t = $p.Range.Text
        if ($found1 -eq $null -and $t -like "to design the output*") {
            $found1 = $p
        } elseif ($found2 -eq $null -and $t -like "May need a bandpass*") {
            $found2 = $p
        } elseif ($found3 -eq $null -and $t -like "Think how*") {
            $found3 = $p
            break
        }
    }
}

if ($found1 -eq $null -or $found2 -eq $null -or $found3 -eq $null) {
    throw "Could not re-locate all three inserted paragraphs"
}

$found1.Range.Font.Size = 11
$found1.Range.Font.SizeBi = 11

$found2.Range.Font.Size = 11
$found2.Range.Font.SizeBi = 11

Write-Output "Replaced paragraph with new content and fixed sizes."
